# The table-row repeater tag currently reads "{%tr for row in items %}"
# (split, by coincidence, across two runs as "...in i" + "tems %}").
# We need it to read "{%tr for row in products %}" while keeping the
# literal text around the changed word in its own, separately-formatted
# run (matching how Word naturally re-flows runs on an in-place edit).

$d = $word.ActiveDocument

$r = $d.Content
$r.Find.Execute("items") | Out-Null

# Briefly diverge the formatting of the matched word so Word is forced to
# split it out of its neighbouring runs instead of silently folding the
# replacement back into one contiguous run, then restore the formatting
# so the final run properties are identical to their neighbours again.
$r.Bold = 1
$r.Text = "products"
$r.Bold = 0
